$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 216, shifting existing rows 216:225 down to 217:226.
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row 216 with the latest weekly price record.
$ws.Range("A216").Value = 3
$ws.Range("B216").Value = "Femacal de La Calera"
$ws.Range("C216").Value = "Coquimbo"
$ws.Range("D216").Value = 44509
$ws.Range("E216").Value = 5
$ws.Range("F216").Value = 100112009
$ws.Range("G216").Value = "Acelga"
$ws.Range("H216").Value = "Sin especificar"
$ws.Range("I216").Value = "Primera"
$ws.Range("J216").Value = 250
$ws.Range("K216").Value = 2000
$ws.Range("L216").Value = 2300
$ws.Range("M216").Value = 2144
$ws.Range("N216").Value = "$/docena de atados (6 kilos)"
$ws.Range("O216").Value = "Provincia de Quillota"
$ws.Range("P216").Value = 357
$ws.Range("Q216").Value = 6
$ws.Range("R216").Value = "Hortaliza"
